$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing last row value (D68: 0.2 -> 0.1)
$ws.Range("D68").Value = 0.1

# Add new row 69 with the new monthly data point.
# Force column A to be treated as text (not auto-converted to a date
# serial) while entering the value, then clear the formatting override
# so the cell keeps the workbook's default style, matching the other
# date-label cells in column A (e.g. A68) which carry no explicit style.
$cellA69 = $ws.Cells.Item(69, 1)
$cellA69.NumberFormat = "@"
$cellA69.Value = "01-09-2021"
$cellA69.ClearFormats()

$ws.Range("B69").Value = 0.4
$ws.Range("C69").Value = -0.8
$ws.Range("D69").Value = 0.7
